$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(46, 8).Value = 260463.4
$ws.Cells.Item(46, 9).Value = 200454.25
$ws.Cells.Item(46, 11).Value = 601362.75
$ws.Cells.Item(46, 13).Value = -601243.75
$ws.Cells.Item(60, 8).Value = 260463.4
$ws.Cells.Item(60, 9).Value = 200454.25
$ws.Cells.Item(60, 11).Value = 601362.75
$ws.Cells.Item(60, 13).Value = -600878.75
$ws.Cells.Item(64, 8).Value = 4243.1665
$ws.Cells.Item(64, 9).Value = 2950
$ws.Cells.Item(64, 11).Value = 2950
$ws.Cells.Item(64, 13).Value = -2702
$ws.Cells.Item(67, 8).Value = 4243.1665
$ws.Cells.Item(67, 9).Value = 2950
$ws.Cells.Item(67, 11).Value = 2950
$ws.Cells.Item(67, 13).Value = -2092
$ws.Cells.Item(76, 8).Value = 3508.3333
$ws.Cells.Item(76, 9).Value = 3580
$ws.Cells.Item(76, 10).Value = 3150
$ws.Cells.Item(76, 11).Value = 3580
$ws.Cells.Item(76, 12).Value = 3150
$ws.Cells.Item(76, 13).Value = -3265
$ws.Cells.Item(76, 14).Value = -3780
$ws.Cells.Item(79, 8).Value = 3508.3333
$ws.Cells.Item(79, 9).Value = 3580
$ws.Cells.Item(79, 10).Value = 3150
$ws.Cells.Item(79, 11).Value = 3580
$ws.Cells.Item(79, 12).Value = 3150
$ws.Cells.Item(79, 13).Value = -2488
$ws.Cells.Item(79, 14).Value = -5334
$ws.Cells.Item(112, 8).Value = 43480470
$ws.Cells.Item(112, 9).Value = 500000740
$ws.Cells.Item(112, 10).Value = 2351.9048
$ws.Cells.Item(112, 11).Value = 1500002220
$ws.Cells.Item(112, 12).Value = 7055.714399999999
$ws.Cells.Item(112, 13).Value = -1500001112
$ws.Cells.Item(112, 14).Value = -9271.714399999999
$ws.Cells.Item(116, 8).Value = 1760.091
$ws.Cells.Item(116, 9).Value = 1662.3334
$ws.Cells.Item(116, 11).Value = 1662.3334
$ws.Cells.Item(116, 13).Value = 1779.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 83502056
$ws.Cells.Item(61, 9).Value = 125127400
$ws.Cells.Item(61, 11).Value = 125127400
$ws.Cells.Item(61, 13).Value = -125127188
$ws.Cells.Item(122, 8).Value = 3368980.8
$ws.Cells.Item(122, 9).Value = 1835.9565
$ws.Cells.Item(122, 10).Value = 11113414
$ws.Cells.Item(122, 11).Value = 5507.8695
$ws.Cells.Item(122, 12).Value = 33340242
$ws.Cells.Item(122, 13).Value = -3057.8695
$ws.Cells.Item(122, 14).Value = -33345142
$ws.Cells.Item(136, 8).Value = 83502056
$ws.Cells.Item(136, 9).Value = 125127400
$ws.Cells.Item(136, 11).Value = 375382200
$ws.Cells.Item(136, 13).Value = -375379650

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 33335940
$ws.Cells.Item(105, 9).Value = 50001800
$ws.Cells.Item(105, 10).Value = 4220
$ws.Cells.Item(105, 11).Value = 50001800
$ws.Cells.Item(105, 12).Value = 4220
$ws.Cells.Item(105, 13).Value = -50000053
$ws.Cells.Item(105, 14).Value = -7714

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 3615
$ws.Cells.Item(62, 9).Value = 2902.5
$ws.Cells.Item(62, 10).Value = 3900
$ws.Cells.Item(62, 11).Value = 2902.5
$ws.Cells.Item(62, 12).Value = 3900
$ws.Cells.Item(62, 13).Value = -2278.5
$ws.Cells.Item(62, 14).Value = -5148
$ws.Cells.Item(65, 8).Value = 3615
$ws.Cells.Item(65, 9).Value = 2902.5
$ws.Cells.Item(65, 10).Value = 3900
$ws.Cells.Item(65, 11).Value = 14512.5
$ws.Cells.Item(65, 12).Value = 19500
$ws.Cells.Item(65, 13).Value = -11392.5
$ws.Cells.Item(65, 14).Value = -25740
$ws.Cells.Item(98, 8).Value = 55945
$ws.Cells.Item(98, 10).Value = 55945
$ws.Cells.Item(98, 12).Value = 55945
$ws.Cells.Item(98, 14).Value = -60437
$ws.Cells.Item(99, 8).Value = 1247.5333
$ws.Cells.Item(99, 9).Value = 1000
$ws.Cells.Item(99, 10).Value = 1618.8334
$ws.Cells.Item(99, 11).Value = 1000
$ws.Cells.Item(99, 12).Value = 1618.8334
$ws.Cells.Item(99, 13).Value = 498
$ws.Cells.Item(99, 14).Value = -4614.8334
$ws.Cells.Item(126, 8).Value = 1247.5333
$ws.Cells.Item(126, 9).Value = 1000
$ws.Cells.Item(126, 10).Value = 1618.8334
$ws.Cells.Item(126, 11).Value = 3000
$ws.Cells.Item(126, 12).Value = 4856.5002
$ws.Cells.Item(126, 13).Value = -530
$ws.Cells.Item(126, 14).Value = -9796.5002
$ws.Cells.Item(134, 8).Value = 24770.623
$ws.Cells.Item(134, 9).Value = 686.59375
$ws.Cells.Item(134, 10).Value = 84054.38
$ws.Cells.Item(134, 11).Value = 2059.78125
$ws.Cells.Item(134, 12).Value = 252163.14
$ws.Cells.Item(134, 13).Value = 475.21875
$ws.Cells.Item(134, 14).Value = -257233.14
$ws.Cells.Item(141, 8).Value = 42609.316
$ws.Cells.Item(141, 9).Value = 10296
$ws.Cells.Item(141, 10).Value = 44404.5
$ws.Cells.Item(141, 11).Value = 10296
$ws.Cells.Item(141, 12).Value = 44404.5
$ws.Cells.Item(141, 13).Value = -5116
$ws.Cells.Item(141, 14).Value = -54764.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(40, 8).Value = 258.75
$ws.Cells.Item(40, 9).Value = 46.5
$ws.Cells.Item(40, 10).Value = 471
$ws.Cells.Item(40, 11).Value = 186
$ws.Cells.Item(40, 12).Value = 1884
$ws.Cells.Item(40, 13).Value = -117
$ws.Cells.Item(40, 14).Value = -2022
$ws.Cells.Item(80, 8).Value = 1849.6818
$ws.Cells.Item(80, 9).Value = 1598
$ws.Cells.Item(80, 10).Value = 1923.7059
$ws.Cells.Item(80, 11).Value = 4794
$ws.Cells.Item(80, 12).Value = 5771.1177
$ws.Cells.Item(80, 13).Value = -3858
$ws.Cells.Item(80, 14).Value = -7643.1177
$ws.Cells.Item(81, 8).Value = 17500.334
$ws.Cells.Item(81, 9).Value = 1250.5
$ws.Cells.Item(81, 10).Value = 50000
$ws.Cells.Item(81, 11).Value = 3751.5
$ws.Cells.Item(81, 12).Value = 150000
$ws.Cells.Item(81, 13).Value = -2628.5
$ws.Cells.Item(81, 14).Value = -152246
$ws.Cells.Item(83, 8).Value = 1849.6818
$ws.Cells.Item(83, 9).Value = 1598
$ws.Cells.Item(83, 10).Value = 1923.7059
$ws.Cells.Item(83, 11).Value = 14382
$ws.Cells.Item(83, 12).Value = 17313.3531
$ws.Cells.Item(83, 13).Value = -9702
$ws.Cells.Item(83, 14).Value = -26673.3531
$ws.Cells.Item(84, 8).Value = 17500.334
$ws.Cells.Item(84, 9).Value = 1250.5
$ws.Cells.Item(84, 10).Value = 50000
$ws.Cells.Item(84, 11).Value = 11254.5
$ws.Cells.Item(84, 12).Value = 450000
$ws.Cells.Item(84, 13).Value = -5638.5
$ws.Cells.Item(84, 14).Value = -461232
$ws.Cells.Item(113, 8).Value = 756
$ws.Cells.Item(113, 9).Value = 593.75
$ws.Cells.Item(113, 10).Value = 885.8
$ws.Cells.Item(113, 11).Value = 1781.25
$ws.Cells.Item(113, 12).Value = 2657.4
$ws.Cells.Item(113, 13).Value = 388.75
$ws.Cells.Item(113, 14).Value = -6997.4
$ws.Cells.Item(117, 8).Value = 8334124
$ws.Cells.Item(117, 9).Value = 707.25
$ws.Cells.Item(117, 10).Value = 16667541
$ws.Cells.Item(117, 11).Value = 2121.75
$ws.Cells.Item(117, 12).Value = 50002623
$ws.Cells.Item(117, 13).Value = 1320.25
$ws.Cells.Item(117, 14).Value = -50009507
$ws.Cells.Item(129, 8).Value = 3207271.2
$ws.Cells.Item(129, 10).Value = 4388213
$ws.Cells.Item(129, 12).Value = 13164639
$ws.Cells.Item(129, 14).Value = -13174639
$ws.Cells.Item(131, 8).Value = 820.07574
$ws.Cells.Item(131, 10).Value = 935.7925
$ws.Cells.Item(131, 12).Value = 2807.3775
$ws.Cells.Item(131, 14).Value = -12887.3775

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(4, 8).Value = 5500
$ws.Cells.Item(4, 9).Value = 500
$ws.Cells.Item(4, 10).Value = 7166.6665
$ws.Cells.Item(4, 11).Value = 500
$ws.Cells.Item(4, 12).Value = 7166.6665
$ws.Cells.Item(4, 13).Value = -388
$ws.Cells.Item(4, 14).Value = -7390.6665
$ws.Cells.Item(70, 8).Value = 31275.078
$ws.Cells.Item(70, 9).Value = 43415.5
$ws.Cells.Item(70, 10).Value = 4970.8335
$ws.Cells.Item(70, 11).Value = 43415.5
$ws.Cells.Item(70, 12).Value = 4970.8335
$ws.Cells.Item(70, 13).Value = -43145.5
$ws.Cells.Item(70, 14).Value = -5510.8335
$ws.Cells.Item(73, 8).Value = 31275.078
$ws.Cells.Item(73, 9).Value = 43415.5
$ws.Cells.Item(73, 10).Value = 4970.8335
$ws.Cells.Item(73, 11).Value = 43415.5
$ws.Cells.Item(73, 12).Value = 4970.8335
$ws.Cells.Item(73, 13).Value = -42479.5
$ws.Cells.Item(73, 14).Value = -6842.8335
$ws.Cells.Item(80, 8).Value = 3547.7144
$ws.Cells.Item(80, 9).Value = 2868.6
$ws.Cells.Item(80, 11).Value = 2868.6
$ws.Cells.Item(80, 13).Value = -1870.6
$ws.Cells.Item(83, 8).Value = 3547.7144
$ws.Cells.Item(83, 9).Value = 2868.6
$ws.Cells.Item(83, 11).Value = 14343
$ws.Cells.Item(83, 13).Value = -9351
$ws.Cells.Item(122, 8).Value = 1819.1794
$ws.Cells.Item(122, 9).Value = 1583.5186
$ws.Cells.Item(122, 10).Value = 2349.4167
$ws.Cells.Item(122, 11).Value = 4750.5558
$ws.Cells.Item(122, 12).Value = 7048.250100000001
$ws.Cells.Item(122, 13).Value = -2300.5558
$ws.Cells.Item(122, 14).Value = -11948.2501
$ws.Cells.Item(126, 8).Value = 1840.125
$ws.Cells.Item(126, 9).Value = 1467.8334
$ws.Cells.Item(126, 10).Value = 2957
$ws.Cells.Item(126, 11).Value = 4403.5002
$ws.Cells.Item(126, 12).Value = 8871
$ws.Cells.Item(126, 13).Value = -1933.5002
$ws.Cells.Item(126, 14).Value = -13811
$ws.Cells.Item(134, 8).Value = 32200
$ws.Cells.Item(134, 10).Value = 32200
$ws.Cells.Item(134, 12).Value = 96600
$ws.Cells.Item(134, 14).Value = -101670
$ws.Cells.Item(138, 8).Value = 53185
$ws.Cells.Item(138, 10).Value = 53185
$ws.Cells.Item(138, 12).Value = 53185
$ws.Cells.Item(138, 14).Value = -63465

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 62502228
$ws.Cells.Item(7, 9).Value = 62502228
$ws.Cells.Item(7, 11).Value = 62502228
$ws.Cells.Item(7, 13).Value = -62502116
$ws.Cells.Item(40, 8).Value = 2557.3333
$ws.Cells.Item(40, 9).Value = 2335.3333
$ws.Cells.Item(40, 10).Value = 2668.3333
$ws.Cells.Item(40, 11).Value = 2335.3333
$ws.Cells.Item(40, 12).Value = 2668.3333
$ws.Cells.Item(40, 13).Value = -2199.3333
$ws.Cells.Item(40, 14).Value = -2940.3333
$ws.Cells.Item(82, 8).Value = 1967.8572
$ws.Cells.Item(82, 9).Value = 1048
$ws.Cells.Item(82, 11).Value = 1048
$ws.Cells.Item(82, 13).Value = -687
$ws.Cells.Item(85, 8).Value = 1967.8572
$ws.Cells.Item(85, 9).Value = 1048
$ws.Cells.Item(85, 11).Value = 1048
$ws.Cells.Item(85, 13).Value = 200
$ws.Cells.Item(122, 8).Value = 3431.2917
$ws.Cells.Item(122, 9).Value = 3331.7778
$ws.Cells.Item(122, 11).Value = 9995.3334
$ws.Cells.Item(122, 13).Value = -7545.3334
$ws.Cells.Item(126, 8).Value = 62502228
$ws.Cells.Item(126, 9).Value = 62502228
$ws.Cells.Item(126, 11).Value = 187506684
$ws.Cells.Item(126, 13).Value = -187504214
$ws.Cells.Item(135, 8).Value = 57359.96
$ws.Cells.Item(135, 10).Value = 57359.96
$ws.Cells.Item(135, 12).Value = 57359.96
$ws.Cells.Item(135, 14).Value = -67499.95999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(54, 8).Value = 14599.4
$ws.Cells.Item(54, 10).Value = 14599.4
$ws.Cells.Item(54, 12).Value = 14599.4
$ws.Cells.Item(54, 14).Value = -15639.4
$ws.Cells.Item(81, 8).Value = 2000.1111
$ws.Cells.Item(81, 9).Value = 1500.1666
$ws.Cells.Item(81, 11).Value = 3000.3332
$ws.Cells.Item(81, 13).Value = -1939.3332
$ws.Cells.Item(84, 8).Value = 2000.1111
$ws.Cells.Item(84, 9).Value = 1500.1666
$ws.Cells.Item(84, 11).Value = 15001.666
$ws.Cells.Item(84, 13).Value = -9697.666000000001

$wb.Save()
Write-Host "Applied 260 changes across 8 sheets"